# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) are refreshed for nearly every row;
# rows 38/39 also swap identity (MXToken <-> RenderToken reordered).
#
# Several "Price" strings are plain decimal-looking text (e.g. "7.710",
# "33.40", "17.30") that Excel would otherwise silently reinterpret as
# numbers (dropping trailing zeros). NumberFormat "@" (Text) is applied
# first for those cells so the literal string is preserved, matching the
# source data which stores every Price/Volume cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.279.99'
$ws.Range("E2").Value = '  +5.35%  '
$ws.Range("D3").Value = '1.918.05'
$ws.Range("E3").Value = '  +5.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.06'
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5172'
$ws.Range("E7").Value = '  +4.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.06'
$ws.Range("E8").Value = '  +6.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2988'
$ws.Range("E9").Value = '  +6.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06837'
$ws.Range("E10").Value = '  +7.00%  '
$ws.Range("D11").Value = '1.917.85'
$ws.Range("E11").Value = '  +6.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '17.53'
$ws.Range("E12").Value = '  +4.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07336'
$ws.Range("E13").Value = '  +3.15%  '
$ws.Range("E14").Value = '  +6.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.88'
$ws.Range("E15").Value = '  +7.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.926'
$ws.Range("E16").Value = '  +4.52%  '
$ws.Range("D17").Value = '30.278.72'
$ws.Range("E17").Value = '  +5.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008032'
$ws.Range("E18").Value = '  +9.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9993'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.11'
$ws.Range("E20").Value = '  +6.90%  '
$ws.Range("D21").Value = '2.165.57'
$ws.Range("E21").Value = '  +6.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.867'
$ws.Range("E23").Value = '  +5.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.773'
$ws.Range("E24").Value = '  +8.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.202'
$ws.Range("E25").Value = '  +3.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.94'
$ws.Range("E26").Value = '  +2.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '139.61'
$ws.Range("E27").Value = '  +25.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.30'
$ws.Range("E28").Value = '  +7.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.015'
$ws.Range("E29").Value = '  +6.94%  '
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.299'
$ws.Range("E31").Value = '  +2.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08877'
$ws.Range("E32").Value = '  +6.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.032'
$ws.Range("E33").Value = '  +4.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05147'
$ws.Range("E34").Value = '  +3.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.161'
$ws.Range("E35").Value = '  +6.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7229'
$ws.Range("E36").Value = '  +6.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.682'
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.860'
$ws.Range("E38").Value = '  +7.77%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.329'
$ws.Range("E39").Value = '  +8.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9756'
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01705'
$ws.Range("E41").Value = '  +6.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.228'
$ws.Range("E42").Value = '  +4.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4348'
$ws.Range("E43").Value = '  +5.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '106.19'
$ws.Range("E44").Value = '  +5.12%  '
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.710'
$ws.Range("E46").Value = '  +6.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1280'
$ws.Range("E47").Value = '  +4.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05736'
$ws.Range("E48").Value = '  +4.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.561'
$ws.Range("E49").Value = '  +5.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.40'
$ws.Range("E50").Value = '  +6.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3857'
$ws.Range("E51").Value = '  +6.89%  '
